$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4, column B) was empty; set it to the generated ValueSet name
$ws.Range("B4").Value = "SurspecialitetransversaleVs"

# "Date" row (row 8, column B) gets refreshed to the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
